$d = $word.ActiveDocument

# Update the date line at the top of the document. The string is unique in
# the document, so a plain whole-document Find/Replace is unambiguous. This
# runs first so every offset captured afterwards already accounts for any
# shift caused by the (one character longer) replacement text.
$d.Content.Find.Execute("2024-10-14 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-10-15 Tuesday", 2)

# Update the division problems in the table. A number of the old and new
# values repeat elsewhere in the table (e.g. "19÷4=" is the original text of
# two different cells with two different replacements, and one cell's new
# value equals another cell's old value), so a single blanket Find/Replace
# pass is ambiguous / order-dependent and this runtime's range-scoped Find
# can end up editing the wrong occurrence when a duplicate is created mid
# pass. To sidestep all of that, every target cell is first addressed by
# its own Range (captured before any table edits happen, so the offsets are
# unaffected by later edits as long as replacement lengths match), and is
# switched to a short placeholder token unique to that cell. Only once every
# cell holds a unique placeholder (so no two cells can ever collide) are the
# placeholders swapped for their real final values. Each placeholder is kept
# exactly as long as the text it replaces so that no other cell's captured
# offsets are invalidated along the way.

$tbl = $d.Tables.Item(1)

$cellSpecs = @(
    @{ Row = 1;  Col = 1; Old = "49÷5="; New = "34÷5=" },
    @{ Row = 1;  Col = 2; Old = "87÷8="; New = "44÷8=" },
    @{ Row = 1;  Col = 3; Old = "61÷5="; New = "45÷6=" },
    @{ Row = 1;  Col = 4; Old = "37÷2="; New = "91÷4=" },
    @{ Row = 1;  Col = 5; Old = "43÷4="; New = "14÷9=" },

    @{ Row = 5;  Col = 1; Old = "19÷4="; New = "59÷6=" },
    @{ Row = 5;  Col = 2; Old = "14÷3="; New = "37÷5=" },
    @{ Row = 5;  Col = 3; Old = "12÷6="; New = "76÷3=" },
    @{ Row = 5;  Col = 4; Old = "83÷9="; New = "36÷7=" },
    @{ Row = 5;  Col = 5; Old = "16÷2="; New = "81÷5=" },

    @{ Row = 9;  Col = 1; Old = "54÷7="; New = "44÷5=" },
    @{ Row = 9;  Col = 2; Old = "89÷8="; New = "73÷7=" },
    @{ Row = 9;  Col = 3; Old = "81÷4="; New = "99÷8=" },
    @{ Row = 9;  Col = 4; Old = "19÷4="; New = "63÷5=" },
    @{ Row = 9;  Col = 5; Old = "83÷6="; New = "49÷4=" },

    @{ Row = 13; Col = 1; Old = "62÷4="; New = "12÷6=" },
    @{ Row = 13; Col = 2; Old = "47÷8="; New = "39÷3=" },
    @{ Row = 13; Col = 3; Old = "21÷8="; New = "45÷7=" },
    @{ Row = 13; Col = 4; Old = "38÷9="; New = "55÷7=" },
    @{ Row = 13; Col = 5; Old = "45÷9="; New = "28÷4=" },

    @{ Row = 17; Col = 1; Old = "63÷4="; New = "94÷7=" },
    @{ Row = 17; Col = 2; Old = "73÷9="; New = "92÷4=" },
    @{ Row = 17; Col = 3; Old = "26÷9="; New = "80÷6=" },
    @{ Row = 17; Col = 4; Old = "57÷6="; New = "78÷3=" },
    @{ Row = 17; Col = 5; Old = "28÷4="; New = "39÷3=" }
)

# Pass 0: capture every target cell's Range boundaries and assign each one a
# unique, same-length placeholder token before any edits are made.
for ($i = 0; $i -lt $cellSpecs.Count; $i++) {
    $spec = $cellSpecs[$i]
    $cellRange = $tbl.Cell($spec.Row, $spec.Col).Range
    $cellSpecs[$i].Start = $cellRange.Start
    $cellSpecs[$i].End = $cellRange.End
    $cellSpecs[$i].Placeholder = "@{0:D3}=" -f $i
}

# Pass 1: swap each cell's original value for its unique placeholder, using
# the pre-captured offsets (stable because Old/Placeholder are equal length).
for ($i = 0; $i -lt $cellSpecs.Count; $i++) {
    $spec = $cellSpecs[$i]
    $rng = $d.Range($spec.Start, $spec.End)
    $rng.Find.Execute($spec.Old, $true, $false, $false, $false, $false,
                       $true, 1, $false, $spec.Placeholder, 2)
}

# Pass 2: swap each unique placeholder for the real replacement value. Since
# every placeholder is unique document-wide, a plain content-wide
# Find/Replace is unambiguous regardless of ordering or length changes.
for ($i = 0; $i -lt $cellSpecs.Count; $i++) {
    $spec = $cellSpecs[$i]
    $d.Content.Find.Execute($spec.Placeholder, $true, $false, $false, $false, $false,
                             $true, 1, $false, $spec.New, 2)
}
